# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.054.76"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "1.826.17"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.92%  "
$ws.Range("D5").Value = "'311.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").Value = "'0.4239"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("D8").Value = "'0.3671"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("D9").Value = "'0.07220"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "'0.8451"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").Value = "1.824.37"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "'6.659"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'5.296"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "'0.07042"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "'89.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "'0.000008748"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").Value = "27.084.44"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").Value = "'5.134"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").Value = "'10.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "2.050.33"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'1.980"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "'151.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").Value = "'2.259"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.73%  "
$ws.Range("D28").Value = "'18.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "'5.249"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "'116.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "'0.08720"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").Value = "'0.7358"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "'4.428"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "'0.05241"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").Value = "'7.331"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'0.1689"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'0.5071"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "'8.558"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "'1.981"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.55%  "
$ws.Range("D47").Value = "'0.4729"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'105.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "'0.06324"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").Value = "'1.653"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.00%  "
